$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Progress tracking updates (US3 / Database tasks) ---
# "Structure of database" task moved from Planned to Completed
$ws.Range("F15").Value = "C"
$ws.Range("F16").Value = "C"
$ws.Range("H16").Value = 8

# "Backend Web Api (GET method)" tasks moved from Draft to In-Progress
$ws.Range("F18").Value = "P"
$ws.Range("F19").Value = "P"
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 2

# --- Update view / selection to match where the author left off ---
$ws.Range("K20").Select()
